$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 10.25574820157827
$ws.Range("D2").Value = 8.975944711274103
$ws.Range("E2").Value = 13.70060417621826
$ws.Range("F2").Value = 34.4933259223421
$ws.Range("G2").Value = 38.12757742536483
$ws.Range("H2").Value = 16.1460037260745
$ws.Range("I2").Value = 26.73130812648933
$ws.Range("J2").Value = 10.16132451069943
$ws.Range("L2").Value = 10.32893753701934
$ws.Range("N2").Value = 18.99769123883691
$ws.Range("O2").Value = 25.95845761623621

$ws.Range("C3").Value = 10.25419415395192
$ws.Range("D3").Value = 8.98551772318287
$ws.Range("E3").Value = 13.71301761476383
$ws.Range("F3").Value = 34.32215378162775
$ws.Range("G3").Value = 37.69371076815724
$ws.Range("H3").Value = 16.13245442039478
$ws.Range("I3").Value = 26.63979269713969
$ws.Range("J3").Value = 10.17838943921167
$ws.Range("L3").Value = 10.34618474666646
$ws.Range("N3").Value = 18.40031508502705
$ws.Range("O3").Value = 25.85418634298525

$ws.Range("C4").Value = 10.25496736024115
$ws.Range("D4").Value = 8.992546182110372
$ws.Range("E4").Value = 13.72265027229677
$ws.Range("F4").Value = 34.22483470244289
$ws.Range("G4").Value = 37.43552600402946
$ws.Range("H4").Value = 16.12714361120193
$ws.Range("I4").Value = 26.58940240185017
$ws.Range("J4").Value = 10.18997233868124
$ws.Range("L4").Value = 10.35787407113777
$ws.Range("N4").Value = 18.02485520896364
$ws.Range("O4").Value = 25.79586358367011

$ws.Range("C5").Value = 10.25571774317531
$ws.Range("D5").Value = 8.995699989109644
$ws.Range("E5").Value = 13.72708169799827
$ws.Range("F5").Value = 34.18716484512658
$ws.Range("G5").Value = 37.33250096491668
$ws.Range("H5").Value = 16.12573820316199
$ws.Range("I5").Value = 26.57034141363614
$ws.Range("J5").Value = 10.19497065615716
$ws.Range("L5").Value = 10.3629143955219
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 25.77354764553522

$ws.Range("C6").Value = 10.2558686579164
$ws.Range("D6").Value = 8.996241179479636
$ws.Range("E6").Value = 13.72784810662931
$ws.Range("F6").Value = 34.18103076345408
$ws.Range("G6").Value = 37.31552933438714
$ws.Range("H6").Value = 16.12555072288738
$ws.Range("I6").Value = 26.5672657457999
$ws.Range("J6").Value = 10.1958174348461
$ws.Range("L6").Value = 10.36376806775663
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 25.76993022197651

$ws.Range("C7").Value = 10.25497571658037
$ws.Range("D7").Value = 8.992587542251956
$ws.Range("E7").Value = 13.72270798663786
$ws.Range("F7").Value = 34.22431858262571
$ws.Range("G7").Value = 37.43412755747885
$ws.Range("H7").Value = 16.12712158243223
$ws.Range("I7").Value = 26.58913935436021
$ws.Range("J7").Value = 10.1900386209079
$ws.Range("L7").Value = 10.35794092537158
$ws.Range("N7").Value = 18.02277304767602
$ws.Range("O7").Value = 25.79555672589981

$ws.Range("C8").Value = 10.25485445262171
$ws.Range("D8").Value = 8.979006887707259
$ws.Range("E8").Value = 13.70446727631461
$ws.Range("F8").Value = 34.43271069899556
$ws.Range("G8").Value = 37.97635093970683
$ws.Range("H8").Value = 16.14070899972144
$ws.Range("I8").Value = 26.69855913196358
$ws.Range("J8").Value = 10.16697936907393
$ws.Range("L8").Value = 10.33465639849372
$ws.Range("N8").Value = 18.79364780656866
$ws.Range("O8").Value = 25.92133296682617

$ws.Range("C9").Value = 10.26827549925819
$ws.Range("D9").Value = 8.961489428108681
$ws.Range("E9").Value = 13.68463057890935
$ws.Range("F9").Value = 34.90157944826341
$ws.Range("G9").Value = 39.09920595849644
$ws.Range("H9").Value = 16.19110787293888
$ws.Range("I9").Value = 26.95842681131283
$ws.Range("J9").Value = 10.13051262704167
$ws.Range("L9").Value = 10.29770340712711
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 26.21234824143748

$ws.Range("C10").Value = 10.28638961953558
$ws.Range("D10").Value = 8.954153750359495
$ws.Range("E10").Value = 13.67973727800897
$ws.Range("F10").Value = 35.28053749958867
$ws.Range("G10").Value = 39.95276507461049
$ws.Range("H10").Value = 16.24243436187514
$ws.Range("I10").Value = 27.17588296160245
$ws.Range("J10").Value = 10.10903585112762
$ws.Range("L10").Value = 10.2758410929544
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 26.45193824389075

$ws.Range("C11").Value = 10.29640291300782
$ws.Range("D11").Value = 8.952012810803081
$ws.Range("E11").Value = 13.67960411940139
$ws.Range("F11").Value = 35.459897779369
$ws.Range("G11").Value = 40.34560918664415
$ws.Range("H11").Value = 16.26884059538284
$ws.Range("I11").Value = 27.28030770749324
$ws.Range("J11").Value = 10.10041518649371
$ws.Range("L11").Value = 10.26703866062771
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 26.56622826055476

$ws.Range("C12").Value = 10.30044764780796
$ws.Range("D12").Value = 8.951373534015241
$ws.Range("E12").Value = 13.67985363855861
$ws.Range("F12").Value = 35.52877180817411
$ws.Range("G12").Value = 40.49488415481617
$ws.Range("H12").Value = 16.27927474049112
$ws.Range("I12").Value = 27.3206175292483
$ws.Range("J12").Value = 10.09731563409771
$ws.Range("L12").Value = 10.26386934274421
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 26.61024105061237

$ws.Range("C13").Value = 10.29956532796522
$ws.Range("D13").Value = 8.951503597632469
$ws.Range("E13").Value = 13.67978657735496
$ws.Range("F13").Value = 35.51389690145547
$ws.Range("G13").Value = 40.46271452596196
$ws.Range("H13").Value = 16.27700831258677
$ws.Range("I13").Value = 27.31190241010112
$ws.Range("J13").Value = 10.09797584964782
$ws.Range("L13").Value = 10.26454462532379
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 26.60072994580684

$ws.Range("C14").Value = 10.29673061937078
$ws.Range("D14").Value = 8.951956784298698
$ws.Range("E14").Value = 13.67961864235957
$ws.Range("F14").Value = 35.46554524596029
$ws.Range("G14").Value = 40.35788064500555
$ws.Range("H14").Value = 16.26969033741641
$ws.Range("I14").Value = 27.2836088072991
$ws.Range("J14").Value = 10.10015688154475
$ws.Range("L14").Value = 10.26677463483672
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 26.56983468329439

$ws.Range("C15").Value = 10.29502715164683
$ws.Range("D15").Value = 8.95225668507454
$ws.Range("E15").Value = 13.67955480695894
$ws.Range("F15").Value = 35.43605126447464
$ws.Range("G15").Value = 40.29372954679263
$ws.Range("H15").Value = 16.26526431818738
$ws.Range("I15").Value = 27.26637721108154
$ws.Range("J15").Value = 10.10151429156747
$ws.Range("L15").Value = 10.2681619229725
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 26.55100511959452

$ws.Range("C16").Value = 10.28577073563207
$ws.Range("D16").Value = 8.954317688712706
$ws.Range("E16").Value = 13.67978801038091
$ws.Range("F16").Value = 35.26895207073295
$ws.Range("G16").Value = 39.9271713231511
$ws.Range("H16").Value = 16.24076976162575
$ws.Range("I16").Value = 27.16916719476818
$ws.Range("J16").Value = 10.10962232749499
$ws.Range("L16").Value = 10.27643931841222
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 26.44457336149639

$ws.Range("C17").Value = 10.28054495942445
$ws.Range("D17").Value = 8.955888023426642
$ws.Range("E17").Value = 13.68046642217093
$ws.Range("F17").Value = 35.16819373417128
$ws.Range("G17").Value = 39.7033654192106
$ws.Range("H17").Value = 16.22652284865344
$ws.Range("I17").Value = 27.11092414121871
$ws.Range("J17").Value = 10.11489044382608
$ws.Range("L17").Value = 10.28180968445781
$ws.Range("N17").Value = 20.9381421901699
$ws.Range("O17").Value = 26.3806184330623

$ws.Range("C18").Value = 10.27770625809705
$ws.Range("D18").Value = 8.956903856558442
$ws.Range("E18").Value = 13.68105368136502
$ws.Range("F18").Value = 35.11089929227093
$ws.Range("G18").Value = 39.5750754300194
$ws.Range("H18").Value = 16.21861652075263
$ws.Range("I18").Value = 27.077943807197
$ws.Range("J18").Value = 10.11802871762254
$ws.Range("L18").Value = 10.28500616954736
$ws.Range("N18").Value = 20.79000725568364
$ws.Range("O18").Value = 26.34433383519426

$ws.Range("C19").Value = 10.27677386884377
$ws.Range("D19").Value = 8.957267158325308
$ws.Range("E19").Value = 13.68128639011132
$ws.Range("F19").Value = 35.09161501483199
$ws.Range("G19").Value = 39.53171807743821
$ws.Range("H19").Value = 16.21598920358623
$ws.Range("I19").Value = 27.06686719854781
$ws.Range("J19").Value = 10.1191098758988
$ws.Range("L19").Value = 10.28610693461117
$ws.Range("N19").Value = 20.73962067985786
$ws.Range("O19").Value = 26.33213531316704

$ws.Range("C20").Value = 10.28108397978834
$ws.Range("D20").Value = 8.955709206869749
$ws.Range("E20").Value = 13.68037381653585
$ws.Range("F20").Value = 35.17885176091223
$ws.Range("G20").Value = 39.72714571648023
$ws.Range("H20").Value = 16.22800967227056
$ws.Range("I20").Value = 27.11707062021442
$ws.Range("J20").Value = 10.11431844882357
$ws.Range("L20").Value = 10.28122686701995
$ws.Range("N20").Value = 20.96544799484618
$ws.Range("O20").Value = 26.38737493709641

$ws.Range("C21").Value = 10.29755639520365
$ws.Range("D21").Value = 8.951819023872796
$ws.Range("E21").Value = 13.67965983716093
$ws.Range("F21").Value = 35.47972179541733
$ws.Range("G21").Value = 40.38866008521513
$ws.Range("H21").Value = 16.27182804914001
$ws.Range("I21").Value = 27.29189872801589
$ws.Range("J21").Value = 10.09951178696438
$ws.Range("L21").Value = 10.26611517981393
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("O21").Value = 26.57888969274275

$ws.Range("C22").Value = 10.30979517933281
$ws.Range("D22").Value = 8.95027562090881
$ws.Range("E22").Value = 13.68094095816693
$ws.Range("F22").Value = 35.68189594445123
$ws.Range("G22").Value = 40.82392603959535
$ws.Range("H22").Value = 16.30299689102792
$ws.Range("I22").Value = 27.41061395752932
$ws.Range("J22").Value = 10.09079580500545
$ws.Range("L22").Value = 10.2571943953773
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 26.70831781823653

$ws.Range("C23").Value = 10.30312904882558
$ws.Range("D23").Value = 8.951008143760589
$ws.Range("E23").Value = 13.68009764378369
$ws.Range("F23").Value = 35.57350140083808
$ws.Range("G23").Value = 40.59139490749993
$ws.Range("H23").Value = 16.28613167351545
$ws.Range("I23").Value = 27.3468543214083
$ws.Range("J23").Value = 10.09535987393048
$ws.Range("L23").Value = 10.26186827498677
$ws.Range("N23").Value = 21.92877110912574
$ws.Range("O23").Value = 26.63885931430897

$ws.Range("C24").Value = 10.28083977250388
$ws.Range("D24").Value = 8.9557896977356
$ws.Range("E24").Value = 13.68041506911904
$ws.Range("F24").Value = 35.17403129195932
$ws.Range("G24").Value = 39.71639345259885
$ws.Range("H24").Value = 16.22733659295999
$ws.Range("I24").Value = 27.11429022403617
$ws.Range("J24").Value = 10.11457670645486
$ws.Range("L24").Value = 10.28149001926079
$ws.Range("N24").Value = 20.95310750188672
$ws.Range("O24").Value = 26.38431881195814

$ws.Range("C25").Value = 10.26318924705725
$ws.Range("D25").Value = 8.965254705996339
$ws.Range("E25").Value = 13.68829407869677
$ws.Range("F25").Value = 34.76851179398107
$ws.Range("G25").Value = 38.78985591801244
$ws.Range("H25").Value = 16.17494762661981
$ws.Range("I25").Value = 26.88338189682043
$ws.Range("J25").Value = 10.13944290318591
$ws.Range("L25").Value = 10.30677011373446
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 26.1289950923724

